# preparation publication 0.2.0
# - bump Version to 0.2.0
# - bump Date to the new publication timestamp
# - add a new "Jurisdiction" / "iso:code:3166:FR" row to the Metadata sheet
#   (pushes Description/Purpose/Copyright/Source/Target down by one row)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- bump Version & Date -------------------------------------------------
$ws.Range("B3").Value = "0.2.0"
$ws.Range("B8").Value = "2023-10-20T08:59:58+00:00"

# --- make room for the new row: stamp the data-row style onto row 16 -----
# (row 16 does not exist yet; copy row 15's formatting down first so the
#  freshly created cells share the same style as the rest of the table)
$ws.Cells.Item(15, 1).Copy($ws.Cells.Item(16, 1))
$ws.Cells.Item(15, 2).Copy($ws.Cells.Item(16, 2))

# --- shift rows 11-15 down to 12-16 (bottom-up) ---------------------------
$ws.Range("A16").Value = $ws.Range("A15").Text
$ws.Range("B16").Value = $ws.Range("B15").Text

$ws.Range("A15").Value = $ws.Range("A14").Text
$ws.Range("B15").Value = $ws.Range("B14").Text

$ws.Range("A14").Value = $ws.Range("A13").Text
$ws.Range("B14").Value = $ws.Range("B13").Text

$ws.Range("A13").Value = $ws.Range("A12").Text
$ws.Range("B13").Value = $ws.Range("B12").Text

$ws.Range("A12").Value = $ws.Range("A11").Text
$ws.Range("B12").Value = $ws.Range("B11").Text

# --- write the new Jurisdiction row at row 11 -----------------------------
$ws.Range("A11").Value = "Jurisdiction"
$ws.Range("B11").Value = "iso:code:3166:FR"
